# Fix duplicate "ENG" node entry with bad coordinates (grooming algorithm permission
# issue) by removing the erroneous row (ID 22, Node ENG, Location 32.661244, 51.677711).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(23).Delete()
